# Add three new data rows (97, 98, 99) to the active sheet, matching the
# existing table layout (columns A..T).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = $ws.Range("D96").NumberFormat

$newRows = @(
    @(4, "Feria Lagunitas de Puerto Montt", "Los Lagos", 44911, 10, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Brooks",  "Primera", 800, 7000, 7500,  7250, "$/bandeja 10 kilos", "Provincia de Curicó", 725, 10),
    @(4, "Feria Lagunitas de Puerto Montt", "Los Lagos", 44911, 10, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Lapins",  "Primera", 800, 7000, 7500,  7250, "$/bandeja 10 kilos", "Provincia de Curicó", 725, 10),
    @(4, "Feria Lagunitas de Puerto Montt", "Los Lagos", 44911, 10, "Fruta", 100103, "Frutos de hueso (carozo)", 100103001, "Cereza", "Rainier", "Primera", 600, 9000, 10000, 9500, "$/bandeja 10 kilos", "Provincia de Curicó", 950, 10)
)

$startRow = 97

for ($r = 0; $r -lt $newRows.Length; $r++) {
    $rowIndex = $startRow + $r
    $rowValues = $newRows[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($rowIndex, $c + 1).Value = $rowValues[$c]
    }
    # Column D holds a date value; reapply the same date/time number format
    # used by the rest of the column so the cell keeps the correct style.
    $ws.Cells.Item($rowIndex, 4).NumberFormat = $dateFormat
}
